$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: rename "reason" column to "tech_reason"
$ws.Range("C1").Value = "tech_reason"

# Row 2
$ws.Range("C2").Value = "The candidate lacks proficiency in CSS and NextJS which are essential for the job. However, the candidate has strong skills in ReactJS, JavaScript, and Frontend Development which are key requirements for the position."
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = "The candidate has a strong interest in AI and Japanese work culture. They have expressed a willingness to adapt to new cultures and have highlighted their strengths in teamwork and problem-solving. These traits align well with the company's soft skill requirements."

# Row 3
$ws.Range("C3").Value = "The candidate has experience in relevant technologies like ReactJS and has worked on projects involving frontend development. However, the candidate lacks experience with NextJS, which is mentioned as a required skill. Due to this, the score is moderate rather than high."
$ws.Range("E3").Value = "The applicant shows strong communication skills and a willingness to adapt to new cultures. They have also demonstrated an ability to work well in a team and handle setbacks effectively. However, the candidate should work on improving time management skills and the ability to acquire new skills more quickly."

# Row 4
$ws.Range("C4").Value = "The candidate lacks the required skills of ReactJS and NextJS which are essential for the job. Although the candidate has experience in Django Rest Framework and ReactJS, the absence of NextJS and the primary lack of ReactJS skills impacts the overall aptness for the job."
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "The candidate has shown strong interest in AI and Japan's work culture. They have displayed adaptability, teamwork, and problem-solving skills in their responses, aligning with the soft skills the company is looking for. Their long-term career plan also indicates commitment."
